$d = $word.ActiveDocument

# Mapping of old text -> new text, applied in document order.
$replacements = @(
    @{ Old = "2025-04-14 Monday"; New = "2025-04-15 Tuesday" },
    @{ Old = "44×37=1628";        New = "16×72=1152" },
    @{ Old = "91×53=4823";        New = "15×62=930" },
    @{ Old = "72×86=6192";        New = "80×90=7200" },
    @{ Old = "54×91=4914";        New = "41×22=902" },
    @{ Old = "31×30=930";         New = "55×65=3575" },
    @{ Old = "53×78=4134";        New = "25×88=2200" },
    @{ Old = "73×89=6497";        New = "14×36=504" },
    @{ Old = "22×34=748";         New = "91×96=8736" },
    @{ Old = "31×64=1984";        New = "50×93=4650" },
    @{ Old = "48×87=4176";        New = "71×68=4828" },
    @{ Old = "78×41=3198";        New = "52×44=2288" },
    @{ Old = "31×14=434";         New = "36×12=432" },
    @{ Old = "64×73=4672";        New = "81×38=3078" },
    @{ Old = "48×95=4560";        New = "19×69=1311" },
    @{ Old = "86×14=1204";        New = "30×39=1170" },
    @{ Old = "20×85=1700";        New = "60×40=2400" },
    @{ Old = "56×95=5320";        New = "13×35=455" },
    @{ Old = "65×40=2600";        New = "99×89=8811" },
    @{ Old = "17×62=1054";        New = "93×46=4278" },
    @{ Old = "95×56=5320";        New = "83×93=7719" },
    @{ Old = "42×27=1134";        New = "42×20=840" },
    @{ Old = "61×50=3050";        New = "33×30=990" },
    @{ Old = "76×51=3876";        New = "80×22=1760" },
    @{ Old = "51×78=3978";        New = "41×21=861" },
    @{ Old = "14×12=168";         New = "23×59=1357" }
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
